# Auto-generated Excel COM-interop script
# Applies market-data refresh values to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3355974.8
$ws.Range("J17").Value = 3635431
$ws.Range("L17").Value = 10906293
$ws.Range("N17").Value = -10906629
$ws.Range("H92").Value = 628.7778
$ws.Range("J92").Value = 510.8
$ws.Range("L92").Value = 510.8
$ws.Range("N92").Value = -3006.8
$ws.Range("H106").Value = 16967.48
$ws.Range("I106").Value = 4012.6875
$ws.Range("K106").Value = 4012.6875
$ws.Range("M106").Value = -3381.6875
$ws.Range("H107").Value = 933.875
$ws.Range("I107").Value = 787.1875
$ws.Range("J107").Value = 1227.25
$ws.Range("K107").Value = 787.1875
$ws.Range("L107").Value = 1227.25
$ws.Range("M107").Value = 1132.8125
$ws.Range("N107").Value = -5067.25
$ws.Range("H116").Value = 6344.8667
$ws.Range("I116").Value = 5043
$ws.Range("K116").Value = 5043
$ws.Range("M116").Value = -1601
$ws.Range("H132").Value = 2688.147
$ws.Range("I132").Value = 2807.862
$ws.Range("J132").Value = 1993.8
$ws.Range("K132").Value = 8423.585999999999
$ws.Range("L132").Value = 5981.4
$ws.Range("M132").Value = -5893.585999999999
$ws.Range("N132").Value = -11041.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4042
$ws.Range("I2").Value = 3114.5667
$ws.Range("K2").Value = 3114.5667
$ws.Range("M2").Value = -3001.5667
$ws.Range("H22").Value = 1499.5
$ws.Range("I22").Value = 2500
$ws.Range("K22").Value = 2500
$ws.Range("M22").Value = -2201
$ws.Range("H32").Value = 19765.715
$ws.Range("I32").Value = 20643.785
$ws.Range("K32").Value = 20643.785
$ws.Range("M32").Value = -20356.785
$ws.Range("H45").Value = 8875.6
$ws.Range("I45").Value = 10365.333
$ws.Range("J45").Value = 2916.6667
$ws.Range("K45").Value = 10365.333
$ws.Range("L45").Value = 2916.6667
$ws.Range("M45").Value = -9988.333000000001
$ws.Range("N45").Value = -3670.6667
$ws.Range("H61").Value = 4116.4707
$ws.Range("I61").Value = 1466.5834
$ws.Range("J61").Value = 10476.2
$ws.Range("K61").Value = 1466.5834
$ws.Range("L61").Value = 10476.2
$ws.Range("M61").Value = -1254.5834
$ws.Range("N61").Value = -10900.2
$ws.Range("H102").Value = 2055
$ws.Range("I102").Value = 1759.3334
$ws.Range("J102").Value = 3533.3333
$ws.Range("K102").Value = 1759.3334
$ws.Range("L102").Value = 3533.3333
$ws.Range("M102").Value = -137.3334
$ws.Range("N102").Value = -6777.3333
$ws.Range("H110").Value = 1299.174
$ws.Range("I110").Value = 1168.55
$ws.Range("J110").Value = 2170
$ws.Range("K110").Value = 1168.55
$ws.Range("L110").Value = 2170
$ws.Range("M110").Value = 876.45
$ws.Range("N110").Value = -6260
$ws.Range("H116").Value = 4042
$ws.Range("I116").Value = 3114.5667
$ws.Range("K116").Value = 3114.5667
$ws.Range("M116").Value = -820.5666999999999
$ws.Range("H136").Value = 4116.4707
$ws.Range("I136").Value = 1466.5834
$ws.Range("J136").Value = 10476.2
$ws.Range("K136").Value = 4399.7502
$ws.Range("L136").Value = 31428.6
$ws.Range("M136").Value = -1849.7502
$ws.Range("N136").Value = -36528.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4042
$ws.Range("I3").Value = 3114.5667
$ws.Range("K3").Value = 3114.5667
$ws.Range("M3").Value = -3000.5667
$ws.Range("H94").Value = 390.8889
$ws.Range("I94").Value = 316.0625
$ws.Range("J94").Value = 989.5
$ws.Range("K94").Value = 316.0625
$ws.Range("L94").Value = 989.5
$ws.Range("M94").Value = 134.9375
$ws.Range("N94").Value = -1891.5
$ws.Range("H105").Value = 3038.1667
$ws.Range("I105").Value = 3332.1765
$ws.Range("J105").Value = 2653.6924
$ws.Range("K105").Value = 3332.1765
$ws.Range("L105").Value = 2653.6924
$ws.Range("M105").Value = -1585.1765
$ws.Range("N105").Value = -6147.6924
$ws.Range("H107").Value = 2449.5264
$ws.Range("I107").Value = 2700.5386
$ws.Range("K107").Value = 2700.5386
$ws.Range("M107").Value = -780.5385999999999
$ws.Range("H134").Value = 7218.788
$ws.Range("I134").Value = 4008.7693
$ws.Range("K134").Value = 12026.3079
$ws.Range("M134").Value = -9491.3079
$ws.Range("H141").Value = 47143.566
$ws.Range("J141").Value = 47143.566
$ws.Range("L141").Value = 47143.566
$ws.Range("N141").Value = -57503.566

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10269.4
$ws.Range("I58").Value = 1839.8
$ws.Range("J58").Value = 18699
$ws.Range("K58").Value = 1839.8
$ws.Range("L58").Value = 18699
$ws.Range("M58").Value = -1636.8
$ws.Range("N58").Value = -19105
$ws.Range("H132").Value = 36471.137
$ws.Range("I132").Value = 27354.8
$ws.Range("J132").Value = 47196.234
$ws.Range("K132").Value = 82064.39999999999
$ws.Range("L132").Value = 141588.702
$ws.Range("M132").Value = -79534.39999999999
$ws.Range("N132").Value = -146648.702
$ws.Range("H136").Value = 10269.4
$ws.Range("I136").Value = 1839.8
$ws.Range("J136").Value = 18699
$ws.Range("K136").Value = 5519.4
$ws.Range("L136").Value = 56097
$ws.Range("M136").Value = -2969.4
$ws.Range("N136").Value = -61197

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1239.3334
$ws.Range("J5").Value = 2053.4707
$ws.Range("L5").Value = 6160.4121
$ws.Range("N5").Value = -6384.4121
$ws.Range("H17").Value = 3778
$ws.Range("I17").Value = 2995
$ws.Range("J17").Value = 4300
$ws.Range("K17").Value = 8985
$ws.Range("L17").Value = 12900
$ws.Range("M17").Value = -8816
$ws.Range("N17").Value = -13238
$ws.Range("H107").Value = 1682.6666
$ws.Range("I107").Value = 2020.2
$ws.Range("K107").Value = 6060.6
$ws.Range("M107").Value = -4140.6
$ws.Range("H135").Value = 1239.3334
$ws.Range("J135").Value = 2053.4707
$ws.Range("L135").Value = 18481.2363
$ws.Range("N135").Value = -23551.2363
$ws.Range("H139").Value = 901
$ws.Range("I139").Value = 901
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 2703
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 2437
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 294941.28
$ws.Range("I113").Value = 294941.28
$ws.Range("K113").Value = 294941.28
$ws.Range("M113").Value = -292771.28
$ws.Range("H126").Value = 2926.5715
$ws.Range("I126").Value = 2947.6667
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 8843.000100000001
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -6373.000100000001
$ws.Range("N126").Value = -13340
$ws.Range("H132").Value = 2457.9285
$ws.Range("I132").Value = 2496.56
$ws.Range("J132").Value = 2136
$ws.Range("K132").Value = 7489.68
$ws.Range("L132").Value = 6408
$ws.Range("M132").Value = -4959.68
$ws.Range("N132").Value = -11468

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3781.7778
$ws.Range("I93").Value = 5621.8887
$ws.Range("K93").Value = 5621.8887
$ws.Range("M93").Value = -4373.8887
$ws.Range("H132").Value = 3654.1538
$ws.Range("I132").Value = 3703
$ws.Range("J132").Value = 3544.25
$ws.Range("K132").Value = 11109
$ws.Range("L132").Value = 10632.75
$ws.Range("M132").Value = -8579
$ws.Range("N132").Value = -15692.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3335515.5
$ws.Range("I136").Value = 4167614.5
$ws.Range("J136").Value = 7119.6665
$ws.Range("K136").Value = 12502843.5
$ws.Range("L136").Value = 21358.9995
$ws.Range("M136").Value = -12500293.5
$ws.Range("N136").Value = -26458.9995
